# Applies the commit's changes:
#  1. Clear cell B5 on "ODI Batting" (the cell is removed entirely).
#  2. Add a new worksheet "ODI Batting Extra" at the end of the workbook,
#     with a header row reusing the existing bold/bordered header style,
#     and four data rows of match-level batting extras.

$wb = $excel.ActiveWorkbook

# --- 1. "ODI Batting": drop the stray empty B5 cell -----------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B5").Value = $null

# --- 2. Add the new "ODI Batting Extra" sheet after the last sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Reuse the header style (bold font, thin border, centered) already used by
# the other sheets' header rows, by copying its formatting over.
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1:D1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

function Set-TextValue($range, $text) {
    # Force the cell to stay text (avoid Excel auto-converting numeric- or
    # percent-looking strings), then drop back to the default style so no
    # extra number-format styles get introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-BlankCell($range) {
    # Leave a present-but-empty cell (matches the source data's blank
    # fields) without forcing any new style into the workbook.
    $range.NumberFormat = "@"
    $range.Value = ""
    $range.Style = "Normal"
}

# Row 2 - match 4564
Set-TextValue $extra.Range("A2") "4564"
$extra.Range("B2").Value = 9
Set-TextValue $extra.Range("C2") "0"
Set-TextValue $extra.Range("D2") "0"
Set-TextValue $extra.Range("E2") "0.96%"
Set-TextValue $extra.Range("F2") "NO"

# Row 3 - match 4565 (batting position / 4s / 6s / % unknown)
Set-TextValue $extra.Range("A3") "4565"
Set-BlankCell $extra.Range("B3")
Set-BlankCell $extra.Range("C3")
Set-BlankCell $extra.Range("D3")
Set-BlankCell $extra.Range("E3")
Set-TextValue $extra.Range("F3") "NO"

# Row 4 - match 4567
Set-TextValue $extra.Range("A4") "4567"
$extra.Range("B4").Value = 10
Set-TextValue $extra.Range("C4") "0"
Set-TextValue $extra.Range("D4") "0"
Set-TextValue $extra.Range("E4") "0.95%"
Set-TextValue $extra.Range("F4") "NO"

# Row 5 - match 4728 (all extra fields unknown)
Set-TextValue $extra.Range("A5") "4728"
Set-BlankCell $extra.Range("B5")
Set-BlankCell $extra.Range("C5")
Set-BlankCell $extra.Range("D5")
Set-BlankCell $extra.Range("E5")
Set-BlankCell $extra.Range("F5")
